$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)

# ------------------------------------------------------------------
# 1) The paragraph used to contain a real Word field:
#       { m:self. }
#    rendered via fldChar begin / instrText / fldChar end runs.
#    The new version wants plain literal text "{m:" + "self" + ".}"
#    (with "self" keeping its special orange colour) instead of an
#    actual field. So: drop the field object itself first.
# ------------------------------------------------------------------
$flds = $hdr.Range.Fields
$fld = $flds.Item(1)
$fld.Delete()

# ------------------------------------------------------------------
# 2) Re-insert the literal text that used to be the field's code,
#    right before the "Invalid query statement" run. Each fragment
#    is (re-)located with a fresh Find so stale Range offsets from a
#    previous InsertBefore never bite us. Inserting back-to-front
#    (the part that must end up right next to "Invalid..." first)
#    keeps the final reading order correct.
# ------------------------------------------------------------------
$rClose = $hdr.Range.Duplicate
$rClose.Find.Execute("Invalid query statement") | Out-Null
$rClose.Collapse(1)
$rClose.InsertBefore(".}")

$rSelf = $hdr.Range.Duplicate
$rSelf.Find.Execute("Invalid query statement") | Out-Null
$rSelf.Collapse(1)
$rSelf.InsertBefore("self")

$rOpen = $hdr.Range.Duplicate
$rOpen.Find.Execute("Invalid query statement") | Out-Null
$rOpen.Collapse(1)
$rOpen.InsertBefore("{m:")

# ------------------------------------------------------------------
# 3) The freshly inserted text inherited the bold/red look of the
#    "Invalid query statement..." run it was merged into. Put it
#    back to plain (non-bold, automatic colour) formatting, except
#    for "self" which keeps the original field's orange colour.
# ------------------------------------------------------------------
$rOpenFmt = $hdr.Range.Duplicate
$rOpenFmt.Find.Execute("{m:") | Out-Null
$rOpenFmt.Font.Bold = 0
$rOpenFmt.Font.Color = -16777216

$rSelfFmt = $hdr.Range.Duplicate
$rSelfFmt.Find.Execute("self") | Out-Null
$rSelfFmt.Font.Bold = 0
$rSelfFmt.Font.Color = 683235

$rCloseFmt = $hdr.Range.Duplicate
$rCloseFmt.Find.Execute(".}") | Out-Null
$rCloseFmt.Font.Bold = 0
$rCloseFmt.Font.Color = -16777216

# ------------------------------------------------------------------
# 4) Prefix the red/bold error message with the "    <---" marker.
# ------------------------------------------------------------------
$rMsg = $hdr.Range.Duplicate
$rMsg.Find.Execute("Invalid query statement") | Out-Null
$rMsg.Collapse(1)
$rMsg.InsertBefore("    <---")
